# Apply the target edit to the "数字激光音视盘机" (digital laser audio/video disc
# player) sheet:
#   1. Within every year group of 4 data rows (A/B/C/D sub-periods), the "B"
#      quarter row and "C" quarter row have their data (columns A:E) swapped
#      while staying on the same row numbers. 2004 has no "A" row, so its
#      group starts directly with the B/C pair on rows 2-3.
#   2. Columns F ("...产销率") and G ("...销售量") are removed entirely, which
#      also drops their header cells (F1/G1) and shrinks the used range from
#      A1:G48 down to A1:E48.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs (1-based) whose A:E contents must be swapped.
$swapPairs = @(
    @(2,3),
    @(6,7),
    @(10,11),
    @(14,15),
    @(18,19),
    @(22,23),
    @(26,27),
    @(30,31),
    @(34,35),
    @(38,39),
    @(42,43),
    @(46,47)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    for ($col = 1; $col -le 5; $col++) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}

# Remove columns F and G (along with their header labels), shifting nothing
# to the right of them since they are the last used columns.
$ws.Range("F1:G48").Delete()

Write-Host "done"
